# Update the "取得日時" (retrieved-at) timestamp on rows 2-6 of the
# "ランサーズ" sheet from 2025-12-07 06:33:09 to 2025-12-07 12:33:26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-07 12:33:26"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
